# --- Commit: Updated symbol list on Wed Jan  4 03:06:45 UTC 2023 with GitHub Actions ---
#
# Automated coinranking.com price-sheet refresh. For every data row (2-51):
#   - column G ("Hora"/hour marker) flips from 2 -> 3
#   - columns D (Price) / E (Volume(1h) % change) get refreshed quotes where the
#     coin still has live market data (rows with "--" placeholders are left alone)
# A handful of low-cap coins (rows 10-15 and 25-26) were re-ranked between snapshots,
# so their Coin name (B) / Link (C) / Price (D) / Volume (E) are replaced wholesale.
#
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($a1, $text) {
    # D/E/G hold numeric- or percent-looking strings (prices, "2.49%", the bare "3")
    # that must stay plain text, exactly as authored upstream (t="inlineStr" in the
    # sheet XML) rather than Excel's normal auto-convert-to-number/percent behaviour.
    # A leading apostrophe is Excel's own quote-prefix for "force this literal as text".
    $ws.Range($a1).Formula = "'" + $text
}

# Row 2
Set-TextCell "D2" '251.72'
Set-TextCell "E2" '2.49%'
Set-TextCell "G2" '3'

# Row 3
Set-TextCell "D3" '28.29'
Set-TextCell "E3" '-3.79%'
Set-TextCell "G3" '3'

# Row 4
Set-TextCell "D4" '5.275'
Set-TextCell "E4" '2.19%'
Set-TextCell "G4" '3'

# Row 5
Set-TextCell "D5" '0.05758'
Set-TextCell "E5" '0.04%'
Set-TextCell "G5" '3'

# Row 6
Set-TextCell "D6" '6.673'
Set-TextCell "E6" '1.50%'
Set-TextCell "G6" '3'

# Row 7
Set-TextCell "E7" '3.03%'
Set-TextCell "G7" '3'

# Row 8
Set-TextCell "D8" '0.8634'
Set-TextCell "E8" '0.72%'
Set-TextCell "G8" '3'

# Row 9
Set-TextCell "D9" '0.9238'
Set-TextCell "E9" '7.44%'
Set-TextCell "G9" '3'

# Row 10 -> WazirX
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextCell "D10" '0.1408'
Set-TextCell "E10" '2.88%'
Set-TextCell "G10" '3'

# Row 11 -> MandalaExchangeToken
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextCell "D11" '0.07150'
Set-TextCell "E11" '1.73%'
Set-TextCell "G11" '3'

# Row 12 -> BitrueCoin
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextCell "D12" '0.03155'
Set-TextCell "E12" '4.37%'
Set-TextCell "G12" '3'

# Row 13 -> BitMartToken
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextCell "D13" '0.09229'
Set-TextCell "E13" '-1.41%'
Set-TextCell "G13" '3'

# Row 14 -> BitForexToken
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextCell "D14" '0.001530'
Set-TextCell "E14" '-0.18%'
Set-TextCell "G14" '3'

# Row 15 -> One
$ws.Range("B15").Value = 'One'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextCell "D15" '0.0006044'
Set-TextCell "E15" '0.20%'
Set-TextCell "G15" '3'

# Row 16
Set-TextCell "D16" '0.005892'
Set-TextCell "E16" '-3.08%'
Set-TextCell "G16" '3'

# Row 17
Set-TextCell "E17" '0.21%'
Set-TextCell "G17" '3'

# Row 18
Set-TextCell "E18" '2.64%'
Set-TextCell "G18" '3'

# Row 19
Set-TextCell "E19" '-2.36%'
Set-TextCell "G19" '3'

# Row 20
Set-TextCell "D20" '0.03405'
Set-TextCell "E20" '3.09%'
Set-TextCell "G20" '3'

# Row 21
Set-TextCell "D21" '0.1314'
Set-TextCell "E21" '2.33%'
Set-TextCell "G21" '3'

# Row 22
Set-TextCell "D22" '3.519'
Set-TextCell "E22" '0.85%'
Set-TextCell "G22" '3'

# Row 23
Set-TextCell "D23" '0.04164'
Set-TextCell "E23" '0.16%'
Set-TextCell "G23" '3'

# Row 24
Set-TextCell "E24" '-1.57%'
Set-TextCell "G24" '3'

# Row 25 -> BitKan
$ws.Range("B25").Value = 'BitKan'
$ws.Range("C25").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-TextCell "D25" '0.001217'
Set-TextCell "E25" '-0.75%'
Set-TextCell "G25" '3'

# Row 26 -> HotbitToken
$ws.Range("B26").Value = 'HotbitToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-TextCell "D26" '0.004924'
Set-TextCell "E26" '19.23%'
Set-TextCell "G26" '3'

# Row 27
Set-TextCell "D27" '0.0001198'
Set-TextCell "E27" '-0.96%'
Set-TextCell "G27" '3'

# Row 28
Set-TextCell "D28" '0.0001936'
Set-TextCell "E28" '33.59%'
Set-TextCell "G28" '3'

# Row 29
Set-TextCell "G29" '3'

# Row 30
Set-TextCell "G30" '3'

# Row 31
Set-TextCell "G31" '3'

# Row 32
Set-TextCell "G32" '3'

# Row 33
Set-TextCell "G33" '3'

# Row 34
Set-TextCell "G34" '3'

# Row 35
Set-TextCell "G35" '3'

# Row 36
Set-TextCell "G36" '3'

# Row 37
Set-TextCell "G37" '3'

# Row 38
Set-TextCell "G38" '3'

# Row 39
Set-TextCell "G39" '3'

# Row 40
Set-TextCell "D40" '0.03842'
Set-TextCell "E40" '3.20%'
Set-TextCell "G40" '3'

# Row 41
Set-TextCell "D41" '0.005682'
Set-TextCell "E41" '-3.87%'
Set-TextCell "G41" '3'

# Row 42
Set-TextCell "D42" '0.1084'
Set-TextCell "E42" '1.23%'
Set-TextCell "G42" '3'

# Row 43
Set-TextCell "D43" '0.002427'
Set-TextCell "E43" '-0.54%'
Set-TextCell "G43" '3'

# Row 44
Set-TextCell "D44" '0.009733'
Set-TextCell "E44" '16.04%'
Set-TextCell "G44" '3'

# Row 45
Set-TextCell "D45" '0.00005267'
Set-TextCell "E45" '-0.16%'
Set-TextCell "G45" '3'

# Row 46
Set-TextCell "D46" '0.00000000749'
Set-TextCell "E46" '-0.12%'
Set-TextCell "G46" '3'

# Row 47
Set-TextCell "D47" '0.08493'
Set-TextCell "E47" '46.37%'
Set-TextCell "G47" '3'

# Row 48
Set-TextCell "D48" '0.002175'
Set-TextCell "E48" '-11.13%'
Set-TextCell "G48" '3'

# Row 49
Set-TextCell "D49" '0.00002098'
Set-TextCell "E49" '-0.12%'
Set-TextCell "G49" '3'

# Row 50
Set-TextCell "D50" '0.0001998'
Set-TextCell "E50" '-0.12%'
Set-TextCell "G50" '3'

# Row 51
Set-TextCell "G51" '3'
